# Auto-generated edit script: updates computed profit/price columns (H-N)
# across multiple sheets, as captured by the source diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3143.4167
$ws.Range("I76").Value = 3107.2856
$ws.Range("J76").Value = 3194
$ws.Range("K76").Value = 3107.2856
$ws.Range("L76").Value = 3194
$ws.Range("M76").Value = -2792.2856
$ws.Range("N76").Value = -3824

$ws.Range("H79").Value = 3143.4167
$ws.Range("I79").Value = 3107.2856
$ws.Range("J79").Value = 3194
$ws.Range("K79").Value = 3107.2856
$ws.Range("L79").Value = 3194
$ws.Range("M79").Value = -2015.2856
$ws.Range("N79").Value = -5378

$ws.Range("H98").Value = 1337.1818
$ws.Range("I98").Value = 1701.2858
$ws.Range("J98").Value = 700
$ws.Range("K98").Value = 1701.2858
$ws.Range("L98").Value = 700
$ws.Range("M98").Value = -203.2858000000001
$ws.Range("N98").Value = -3696

$ws.Range("H122").Value = 1337.1818
$ws.Range("I122").Value = 1701.2858
$ws.Range("J122").Value = 700
$ws.Range("K122").Value = 5103.857400000001
$ws.Range("L122").Value = 2100
$ws.Range("M122").Value = -2653.857400000001
$ws.Range("N122").Value = -7000

$ws.Range("H138").Value = 3854.5225
$ws.Range("I138").Value = 1255.625
$ws.Range("J138").Value = 6230.657
$ws.Range("K138").Value = 3766.875
$ws.Range("L138").Value = 18691.971
$ws.Range("M138").Value = 1373.125
$ws.Range("N138").Value = -28971.971

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 37098.855
$ws.Range("I2").Value = 1214.5294
$ws.Range("J2").Value = 92556.45
$ws.Range("K2").Value = 1214.5294
$ws.Range("L2").Value = 92556.45
$ws.Range("M2").Value = -1101.5294
$ws.Range("N2").Value = -92782.45

$ws.Range("H32").Value = 27707.254
$ws.Range("I32").Value = 5010.2954
$ws.Range("J32").Value = 71127.52
$ws.Range("K32").Value = 5010.2954
$ws.Range("L32").Value = 71127.52
$ws.Range("M32").Value = -4723.2954
$ws.Range("N32").Value = -71701.52

$ws.Range("H61").Value = 2468.1738
$ws.Range("I61").Value = 1386.25
$ws.Range("J61").Value = 3045.2
$ws.Range("K61").Value = 1386.25
$ws.Range("L61").Value = 3045.2
$ws.Range("M61").Value = -1174.25
$ws.Range("N61").Value = -3469.2

$ws.Range("H70").Value = 38800
$ws.Range("J70").Value = 38800
$ws.Range("L70").Value = 38800
$ws.Range("N70").Value = -39340

$ws.Range("H73").Value = 38800
$ws.Range("J73").Value = 38800
$ws.Range("L73").Value = 38800
$ws.Range("N73").Value = -40672

$ws.Range("H74").Value = 2363.2188
$ws.Range("I74").Value = 1822.8235
$ws.Range("J74").Value = 2975.6667
$ws.Range("K74").Value = 1822.8235
$ws.Range("L74").Value = 2975.6667
$ws.Range("M74").Value = -948.8235
$ws.Range("N74").Value = -4723.6667

$ws.Range("H77").Value = 2363.2188
$ws.Range("I77").Value = 1822.8235
$ws.Range("J77").Value = 2975.6667
$ws.Range("K77").Value = 9114.1175
$ws.Range("L77").Value = 14878.3335
$ws.Range("M77").Value = -4746.1175
$ws.Range("N77").Value = -23614.3335

$ws.Range("H102").Value = 44899.74
$ws.Range("I102").Value = 59864.53
$ws.Range("J102").Value = 2499.5
$ws.Range("K102").Value = 59864.53
$ws.Range("L102").Value = 2499.5
$ws.Range("M102").Value = -58242.53
$ws.Range("N102").Value = -5743.5

$ws.Range("H116").Value = 37098.855
$ws.Range("I116").Value = 1214.5294
$ws.Range("J116").Value = 92556.45
$ws.Range("K116").Value = 1214.5294
$ws.Range("L116").Value = 92556.45
$ws.Range("M116").Value = 1079.4706
$ws.Range("N116").Value = -97144.45

$ws.Range("H132").Value = 2354.1
$ws.Range("I132").Value = 2238.2258
$ws.Range("J132").Value = 2753.2222
$ws.Range("K132").Value = 6714.6774
$ws.Range("L132").Value = 8259.6666
$ws.Range("M132").Value = -4184.6774
$ws.Range("N132").Value = -13319.6666

$ws.Range("H136").Value = 2468.1738
$ws.Range("I136").Value = 1386.25
$ws.Range("J136").Value = 3045.2
$ws.Range("K136").Value = 4158.75
$ws.Range("L136").Value = 9135.599999999999
$ws.Range("M136").Value = -1608.75
$ws.Range("N136").Value = -14235.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 37098.855
$ws.Range("I3").Value = 1214.5294
$ws.Range("J3").Value = 92556.45
$ws.Range("K3").Value = 1214.5294
$ws.Range("L3").Value = 92556.45
$ws.Range("M3").Value = -1100.5294
$ws.Range("N3").Value = -92784.45

$ws.Range("H99").Value = 1553.9048
$ws.Range("I99").Value = 1214.2727
$ws.Range("J99").Value = 1927.5
$ws.Range("K99").Value = 1214.2727
$ws.Range("L99").Value = 1927.5
$ws.Range("M99").Value = 283.7273
$ws.Range("N99").Value = -4923.5

$ws.Range("H105").Value = 224071.22
$ws.Range("I105").Value = 168780
$ws.Range("J105").Value = 334653.66
$ws.Range("K105").Value = 168780
$ws.Range("L105").Value = 334653.66
$ws.Range("M105").Value = -167033
$ws.Range("N105").Value = -338147.66

$ws.Range("H134").Value = 2635.0286
$ws.Range("I134").Value = 2533.7334
$ws.Range("J134").Value = 3242.8
$ws.Range("K134").Value = 7601.2002
$ws.Range("L134").Value = 9728.400000000001
$ws.Range("M134").Value = -5066.2002
$ws.Range("N134").Value = -14798.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12625.132
$ws.Range("I31").Value = 15056.386
$ws.Range("J31").Value = 4520.952
$ws.Range("K31").Value = 15056.386
$ws.Range("L31").Value = 4520.952
$ws.Range("M31").Value = -14761.386
$ws.Range("N31").Value = -5110.952

$ws.Range("H34").Value = 12625.132
$ws.Range("I34").Value = 15056.386
$ws.Range("J34").Value = 4520.952
$ws.Range("K34").Value = 15056.386
$ws.Range("L34").Value = 4520.952
$ws.Range("M34").Value = -14854.386
$ws.Range("N34").Value = -4924.952

$ws.Range("H58").Value = 8610.941999999999
$ws.Range("I58").Value = 1269.9615
$ws.Range("J58").Value = 29818.223
$ws.Range("K58").Value = 1269.9615
$ws.Range("L58").Value = 29818.223
$ws.Range("M58").Value = -1066.9615
$ws.Range("N58").Value = -30224.223

$ws.Range("H62").Value = 6946619
$ws.Range("I62").Value = 55555556
$ws.Range("J62").Value = 2485.4285
$ws.Range("K62").Value = 55555556
$ws.Range("L62").Value = 2485.4285
$ws.Range("M62").Value = -55554932
$ws.Range("N62").Value = -3733.4285

$ws.Range("H65").Value = 6946619
$ws.Range("I65").Value = 55555556
$ws.Range("J65").Value = 2485.4285
$ws.Range("K65").Value = 277777780
$ws.Range("L65").Value = 12427.1425
$ws.Range("M65").Value = -277774660
$ws.Range("N65").Value = -18667.1425

$ws.Range("H106").Value = 25694
$ws.Range("J106").Value = 25694
$ws.Range("L106").Value = 25694
$ws.Range("N106").Value = -28218

$ws.Range("H132").Value = 2914.6428
$ws.Range("I132").Value = 2199.0476
$ws.Range("K132").Value = 6597.1428
$ws.Range("M132").Value = -4067.1428

$ws.Range("H134").Value = 1095.6
$ws.Range("I134").Value = 1095.6
$ws.Range("K134").Value = 3286.8
$ws.Range("M134").Value = -751.7999999999997

$ws.Range("H136").Value = 8610.941999999999
$ws.Range("I136").Value = 1269.9615
$ws.Range("J136").Value = 29818.223
$ws.Range("K136").Value = 3809.8845
$ws.Range("L136").Value = 89454.66900000001
$ws.Range("M136").Value = -1259.8845
$ws.Range("N136").Value = -94554.66900000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()

$ws.Range("H38").Value = 60.875
$ws.Range("I38").Value = 59.42857
$ws.Range("J38").Value = 62
$ws.Range("K38").Value = 178.28571
$ws.Range("L38").Value = 186
$ws.Range("M38").Value = 168.71429
$ws.Range("N38").Value = -880

$ws.Range("H113").Value = 635.381
$ws.Range("I113").Value = 574.6923
$ws.Range("J113").Value = 734
$ws.Range("K113").Value = 1724.0769
$ws.Range("L113").Value = 2202
$ws.Range("M113").Value = 445.9231
$ws.Range("N113").Value = -6542

$ws.Range("H131").Value = 822.65
$ws.Range("I131").Value = 354.33334
$ws.Range("J131").Value = 852.54254
$ws.Range("K131").Value = 1063.00002
$ws.Range("L131").Value = 2557.62762
$ws.Range("M131").Value = 3976.99998
$ws.Range("N131").Value = -12637.62762

$ws.Range("H133").Value = 7706
$ws.Range("I133").Value = 765
$ws.Range("J133").Value = 12333.333
$ws.Range("K133").Value = 2295
$ws.Range("L133").Value = 36999.999
$ws.Range("M133").Value = 2765
$ws.Range("N133").Value = -47119.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 63584.793
$ws.Range("I70").Value = 109680.84
$ws.Range("J70").Value = 5196.467
$ws.Range("K70").Value = 109680.84
$ws.Range("L70").Value = 5196.467
$ws.Range("M70").Value = -109410.84
$ws.Range("N70").Value = -5736.467

$ws.Range("H73").Value = 63584.793
$ws.Range("I73").Value = 109680.84
$ws.Range("J73").Value = 5196.467
$ws.Range("K73").Value = 109680.84
$ws.Range("L73").Value = 5196.467
$ws.Range("M73").Value = -108744.84
$ws.Range("N73").Value = -7068.467

$ws.Range("H107").Value = 631606.3
$ws.Range("I107").Value = 262.9
$ws.Range("J107").Value = 1683845.4
$ws.Range("K107").Value = 262.9
$ws.Range("L107").Value = 1683845.4
$ws.Range("M107").Value = 1657.1
$ws.Range("N107").Value = -1687685.4

$ws.Range("H122").Value = 3032.24
$ws.Range("I122").Value = 4035.5833
$ws.Range("J122").Value = 2106.077
$ws.Range("K122").Value = 12106.7499
$ws.Range("L122").Value = 6318.231000000001
$ws.Range("M122").Value = -9656.749899999999
$ws.Range("N122").Value = -11218.231

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1622.9032
$ws.Range("I136").Value = 575.63635
$ws.Range("J136").Value = 2198.9
$ws.Range("K136").Value = 1726.90905
$ws.Range("L136").Value = 6596.700000000001
$ws.Range("M136").Value = 823.09095
$ws.Range("N136").Value = -11696.7

